$d = $word.ActiveDocument

# Locate the paragraph that holds the "Ver no Jupiter ..." footer line.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter*Salvar em docx*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Remove the blank paragraph right before it, the "Ver no Jupiter..."
    # paragraph itself, and the following copyright paragraph - i.e. the
    # whole block that made up the trailing site-footer text.
    $startPara = $d.Paragraphs.Item($targetIndex - 1)
    $endPara   = $d.Paragraphs.Item($targetIndex + 1)

    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
